$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, centered, bordered) from H1 into the new I1/J1 headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the I0 (I) and IF (J) data columns for rows 2-69
$data = @{
    2 = @(7, 8)
    3 = @(7, 7)
    4 = @(7, 7)
    5 = @(6, 7)
    6 = @(4, 5)
    7 = @(8, 8)
    8 = @(7, 7)
    9 = @(6, 6)
    10 = @(7, 8)
    11 = @(6, 7)
    12 = @(8, 8)
    13 = @(8, 8)
    14 = @(6, 7)
    15 = @(8, 8)
    16 = @(6, 7)
    17 = @(8, 8)
    18 = @(9, 9)
    19 = @(6, 7)
    20 = @(10, 10)
    21 = @(6, 6)
    22 = @(9, 9)
    23 = @(7, 8)
    24 = @(8, 8)
    25 = @(8, 8)
    26 = @(5, 6)
    27 = @(7, 8)
    28 = @(7, 7)
    29 = @(10, 10)
    30 = @(8, 8)
    31 = @(9, 9)
    32 = @(6, 7)
    33 = @(8, 8)
    34 = @(9, 9)
    35 = @(6, 6)
    36 = @(7, 8)
    37 = @(8, 8)
    38 = @(9, 9)
    39 = @(6, 6)
    40 = @(9, 9)
    41 = @(4, 5)
    42 = @(6, 6)
    43 = @(9, 9)
    44 = @(6, 6)
    45 = @(6, 7)
    46 = @(8, 8)
    47 = @(9, 9)
    48 = @(8, 8)
    49 = @(8, 8)
    50 = @(9, 9)
    51 = @(9, 9)
    52 = @(7, 8)
    53 = @(6, 6)
    54 = @(9, 9)
    55 = @(9, 9)
    56 = @(5, 6)
    57 = @(8, 8)
    58 = @(9, 9)
    59 = @(6, 6)
    60 = @(8, 9)
    61 = @(8, 9)
    62 = @(6, 6)
    63 = @(9, 9)
    64 = @(8, 8)
    65 = @(7, 8)
    66 = @(9, 9)
    67 = @(6, 6)
    68 = @(3, 3)
    69 = @(6, 6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
